$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper note on this engine's InsertXML quirk:
# Range.InsertXML() only cleanly replaces content (keeping the paragraph's
# own <w:pPr>/rsid* attributes intact) when the Range spans the *entire*
# paragraph, including its trailing paragraph mark. If the Range stops
# short of the paragraph mark (e.g. a bookmark or other content follows
# inside the same paragraph), the paragraph wrapper/properties get dropped
# and the new content is appended at the end of the paragraph instead of
# being placed at the selection. So below we always select whole paragraphs
# (via the Paragraphs collection, which includes the trailing mark) and
# supply a full replacement <w:p>...</w:p> including <w:pPr> and the
# original rsid* attributes.
# ---------------------------------------------------------------------------

function Get-ParagraphContainingText($doc, [string]$needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Change 1: "Iteration Duration: 21 November 2016 to 27 November 2016"
#        -> "Iteration Duration: 28 November 2016 to 4 December 2016"
# The author typed over "1" (of "21") with "8", "27" with "4", and the
# second "November" with "December", leaving the cursor (and the _GoBack
# bookmark marking the last edit location) right after "December". We
# reproduce the exact run layout so the runs stay split exactly like a
# real Word editing session would leave them.
# ---------------------------------------------------------------------------

$rPr1 = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/></w:rPr>'
$pPr1 = '<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/></w:rPr></w:pPr>'

$runs1 = ''
$runs1 += '<w:r>' + $rPr1 + '<w:t>Iteration Duration: 2</w:t></w:r>'
$runs1 += '<w:r>' + $rPr1 + '<w:t>8</w:t></w:r>'
$runs1 += '<w:r>' + $rPr1 + '<w:t xml:space="preserve"> November 2016 to </w:t></w:r>'
$runs1 += '<w:r>' + $rPr1 + '<w:t>4</w:t></w:r>'
$runs1 += '<w:r>' + $rPr1 + '<w:t xml:space="preserve"> </w:t></w:r>'
$runs1 += '<w:r>' + $rPr1 + '<w:t>December</w:t></w:r>'
$runs1 += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$runs1 += '<w:r>' + $rPr1 + '<w:t xml:space="preserve"> 2016</w:t></w:r>'

$para1Xml = '<w:p w:rsidR="001F1008" w:rsidRDefault="001F1008" w:rsidP="001F1008">' + $pPr1 + $runs1 + '</w:p>'

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $para1Xml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$p1 = Get-ParagraphContainingText $d "Iteration Duration: 21 November 2016 to 27 November 2016"
if ($null -eq $p1) {
    throw "Could not find the Iteration Duration paragraph to replace"
}
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
[void]$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: drop the old _GoBack bookmark that used to sit right after
# "...the user needs are not clearly defined." (it moved to the location
# above instead). Everything else about that paragraph stays the same.
# ---------------------------------------------------------------------------

$text2 = "-Customer representative unclear about user requirements, the user needs are not clearly defined."
$rPr2 = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr>'
$pPr2 = '<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr>'
$runs2 = '<w:r>' + $rPr2 + '<w:t>' + $text2 + '</w:t></w:r>'
$para2Xml = '<w:p w:rsidR="00934096" w:rsidRDefault="00934096" w:rsidP="00934096">' + $pPr2 + $runs2 + '</w:p>'

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $para2Xml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$p2 = Get-ParagraphContainingText $d "user needs are not clearly defined"
if ($null -eq $p2) {
    throw "Could not find the customer representative paragraph"
}
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
[void]$r2.InsertXML($xml2)

Write-Output "Edit complete"
